# alterando feature selecionar os serviçoes web
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear all manual formatting from the data rows (2-17) - they revert to
# the workbook's default style.
$ws.Range("A2:E17").ClearFormats()

# Re-apply a header style to row 1, but without the blue fill (new cellXfs
# entry: same font/border/alignment as before, fillId changed to "None").
$headerRange = $ws.Range("A1:E1")
$headerRange.Font.Bold = $true
$headerRange.Interior.Pattern = -4142
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Add the new row of data.
$ws.Range("A19").Value = "ejoajiajoejiaoejiaojeoia"
$ws.Range("B19").Value = "FACEPE 59/2022"
$ws.Range("C19").Value = "Thiago Borges Miranda"
$ws.Range("D19").Value = "506070"
$ws.Range("E19").Value = "Dep Fisica"

$ws.Range("A1").Select()
